# Translated guide reshape to en
# The Swedish "ifstatements" guide row had its URL placeholder replaced with a
# plain "(ifstatements_sv.html)" label, and the whole lookup table was then
# re-sorted (descending by "Sv hits april", with "En hits april" as a
# secondary descending tie-breaker).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell that used to hold the ifstatements_sv.html address.
$ws.Range("A20").Value = "(ifstatements_sv.html)"

# Re-sort the whole table (including header) by column B descending,
# breaking ties by column C descending.
$rng = $ws.Range("A1:D45")
$key1 = $ws.Range("B1:B45")
$key2 = $ws.Range("C1:C45")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key1, 0, 2) | Out-Null
$ws.Sort.SortFields.Add($key2, 0, 2) | Out-Null
$ws.Sort.SetRange($rng)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# Restore the view: scroll down and select C53 (mirrors the saved view state).
$ws.Range("A42").Select()
$ws.Range("C53").Select()
